$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("genomic_targeting")

$ws.Range("A16").Value = "GJB2"
$ws.Range("B16").Value = $true
$ws.Range("C16").Value = "GJB2-associated hearing loss- complete genomic targeting"

$ws.Range("C16").Select()
